$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 31: replace the leftover "dada" placeholder text with the real 30th load name
$ws.Range("B31").Value = "Carga 30"

# Row 2: clear the stray wattage value and reset the hourly usage flags back to 0
$ws.Range("C2").ClearContents()
$ws.Range("I2:L2").Value = 0
$ws.Range("Q2:U2").Value = 0

# Restore the previously active selection
$ws.Range("D9").Select()
